$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows 8-10 (old data no longer present after the edit)
$ws.Rows("8:10").Delete()

# New data values for rows 2-7 (columns A,B,C,E,F,G ; D stays blank)
$data = @(
    @("MockHeadersApi", "MockHeadersController", "Get", "A", "Headers", "header*hdr1=ABC&header*hdr2=DEF"),
    @("MockHeadersApi", "MockHeadersController", "Get", "A", "Expected", '[{"Key":"Host","Value":"localhost"},{"Key":"hdr1","Value":"ABC"},{"Key":"hdr2","Value":"DEF"},{"Key":"X-User","Value":"moe@stooges.org"},{"Key":"X-Role","Value":"admin"},{"Key":"X-Role","Value":"user"}]'),
    @("MockHeadersApi", "MockHeadersController", "Get", "B", "Headers", "header*X-User=jill&header*X-Role=user"),
    @("MockHeadersApi", "MockHeadersController", "Get", "B", "Expected", '[{"Key":"Host","Value":"localhost"},{"Key":"X-User","Value":"larry@stooges.org"},{"Key":"X-Role","Value":"admin"},{"Key":"X-Role","Value":"user"}]'),
    @("MockHeadersApi", "MockHeadersController", "Get", "C", "Headers", "header*X-User=jill&header*X-Role=user"),
    @("MockHeadersApi", "MockHeadersController", "Get", "C", "Expected", '[{"Key":"Host","Value":"localhost"},{"Key":"X-User","Value":"curly@stooges.org"},{"Key":"X-Role","Value":"readonly"}]')
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $ws.Cells.Item($row, 6).Value = $entry[4]
    $ws.Cells.Item($row, 7).Value = $entry[5]
    $row++
}

$ws.Range("G8").Select()
